$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive opening text (robust to index shifts).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Grazas por participar*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$r.MoveEnd(1, -1)
$r.Delete()

$p2 = $d.Paragraphs.Item($targetIndex)
$r2 = $p2.Range
$r2.Collapse(1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="37B3C2F6" w14:textId="5EB4D08C" w:rsidR="00134AE6" w:rsidRPr="00F8418D" w:rsidRDefault="00763C02" w:rsidP="004615A9"><w:pPr><w:pStyle w:val="BasicParagraph"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="-72"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Optima-Regular"/><w:sz w:val="20"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:t>Grazas por participar nesta campaña global de medida da contaminación lumínica mediante a observación das estrelas máis febles que podes albiscar. Localizando e observando a  Constelación de Hércules e comparándoa co que aparece nos mapas estelares recollidos neste documento podes saber canto contribúen á contaminación lumínica os sistemas de iluminación que hai no teu barrio ou vila. As túas achegas á base de datos en liña de GLOBE at Night (O MUNDO á Noite) servirán para documentar a calidade do ceo nocturno.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml)
